$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.137.73"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.749.65"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'236.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.5294"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("D8").Value = "'0.2805"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "'0.06180"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "1.747.17"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").Value = "'0.07171"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "'15.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "'0.6454"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'4.630"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("D15").Value = "'78.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "'1.0000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "26.029.66"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'11.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").Value = "'0.000006771"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "1.971.23"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "'4.330"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.21%  "
$ws.Range("D23").Value = "'8.726"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "'5.230"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "'139.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "'1.523"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").Value = "'15.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'1.803"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'3.763"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").Value = "'3.649"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.98%  "
$ws.Range("D33").Value = "'0.04641"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.72%  "
$ws.Range("D34").Value = "'2.647"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'1.007"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").Value = "'0.6327"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("D37").Value = "'2.708"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Value = "'0.01626"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("D39").Value = "'1.975"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").Value = "'0.9998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").Value = "'102.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "'0.3925"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").Value = "'0.7549"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.52%  "
$ws.Range("D44").Value = "'5.073"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D45").Value = "'0.1153"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").Value = "'6.346"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "'0.05341"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").Value = "'54.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.38%  "
$ws.Range("D50").Value = "'0.3474"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").Value = "'7.565"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.38%  "
